# "Add files via upload" -- refresh of the daily IClientBalance export.
# The workbook is re-exported two days later (2024-09-20 -> 2024-09-23), so:
#   1) the sheet name (which embeds the export timestamp) is updated,
#   2) every row's "Dt. Referencia" (column G) date is bumped accordingly,
#   3) a few accounts' "Saldo Previsto" / "Vl. Total" (columns E/H) figures
#      were recalculated upstream and come back with new balances.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet to match the new export run.
$ws.Name = "IClientBalance-20240923-094346-"

# 2) Column G holds the reference date as an Excel serial number
#    (45555 = 2024-09-20, 45558 = 2024-09-23). Update every data row.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
if ($lastRow -lt 274) { $lastRow = 274 }

for ($r = 2; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 7).Value -ne $null) {
        $ws.Cells.Item($r, 7).Value = 45558
    }
}

# 3) Corrected balances for a handful of accounts (row -> new E/H value).
$corrections = @{
    8   = 124.12
    101 = 3781.01
    105 = 13565
    110 = 7665.04
    165 = 68028.48
}

foreach ($row in $corrections.Keys) {
    $value = $corrections[$row]
    $ws.Cells.Item($row, 5).Value = $value   # E: Saldo Previsto
    $ws.Cells.Item($row, 8).Value = $value   # H: Vl. Total
}
